$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Dr. Eman Tantawi, Dr. Majorelle Magdy, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud'
$ws.Range("G3").Value = 'Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Menna tuâ€™Allah Medhat'
$ws.Range("G4").Value = 'Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Rana Abo-Zaid'
$ws.Range("G5").Value = 'Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Nourhan Mahmoud, Dr. Mohammad El-Tanany, Dr. Hanan Ragab, Dr. Nesma, Dr. Veronia Rafat'
$ws.Range("G6").Value = 'Dr. Amira Sobhy, Dr. Gehan Adel, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Nourhan Mahmoud, Dr. Asmaa Reda, Dr. Nahla Nagiub, Dr. Veronia Rafat, Dr. Menna tuâ€™Allah Medhat'
$ws.Range("G7").Value = 'Dr. Gehan Adel, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid'
$ws.Range("G8").Value = 'Administrator, Dr. Manar Montaser, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Shimaa Ahmad Mekki'
$ws.Range("G9").Value = 'Dr. Amira Sobhy, Dr. Gehan Adel, Dr. Manar Montaser, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid'
$ws.Range("G10").Value = 'Dr. Gehan Adel, Dr. Shimaa Ahmad Mekki, Dr. Servinaz Sayed Mohammad, Dr. Heba Mahmoud Ali, Dr. Sara Wael, Dr. Rana Abo-Zaid, Dr. Alshimaa Atef'
$ws.Range("G11").Value = 'Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Asmaa Reda'
$ws.Range("G12").Value = 'Administrator, Dr. Salma El-Gendy'
$ws.Range("G13").Value = 'Dr. Shimaa Ashraf, Dr. Safa Hany, Dr. Mariam Nour El-Din, D Wessam Atef, Dr. Omnia Mohammad'
$ws.Range("G14").Value = 'Dr. Safa Hany, Dr. Shimaa Ashraf'
$ws.Range("G15").Value = 'D Wessam Atef, Dr. Amal Awwad'
$ws.Range("G16").Value = 'Dr. Nourhan Mohammad, Dr. Amal Awwad'
$ws.Range("G17").Value = 'Dr. Marwa Mustafa, Dr. Eman M. Abo-Sakaya, Dr. Nourhan Osama, Dr. Madeha Saeed, Dr. Sarah Abdelmohsen, Dr. Yasmeena Fattoh, Dr. Basma Hamed, Dr. Esraa Mostafa, Dr. Dina Adel, Dr. Arwa Al-Sayed'
$ws.Range("G19").Value = 'D Mariam E. Mohammad, Dr. Sarah Mahdy'
$ws.Range("G22").Value = 'Dr. Nancy Abd Al-Shafy, Dr. Amr Saeed'
$ws.Range("G24").Value = 'Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Youstina Magdy, Dr. Yasmin, Dr. Monica, Dr. Maryam Ashraf, Dr. Aya Emad, Dr. Marina Atef, Dr. Salma Hassan, Dr. Remon, Dr. Neveen Nashaat'
$ws.Range("G25").Value = 'Dr. Ola Abd Al-Fattah, Dr. Youstina Magdy, Dr. Aya Emad, Dr. Abdullah El-Agrody, Dr. Marina Atef, Dr. Remon, Dr. Eman Samir Gabry'
$ws.Range("G27").Value = 'Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Eman Mohammad Al, Dr. Yasmin, Dr. Salma Hassan, Dr. Remon, Dr. Eman Samir Gabry, Dr. Neveen Nashaat'
$ws.Range("G28").Value = 'Dr. Aya Hanafy, Dr. Wafaa Ebida, Dr. Nardine, Dr. Abdullah El-Agrody, Dr. Salma Hassan, Dr. Remon, Dr. Eman Samir Gabry, Dr. Neveen Nashaat'
$ws.Range("G29").Value = 'Dr. Naema Gomaa, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Remon, Dr. Eman Samir Gabry, Dr. Neveen Nashaat'
$ws.Range("G30").Value = 'Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Rana Abo-Zaid'
$ws.Range("G31").Value = 'Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Menna tuâ€™Allah Medhat'
$ws.Range("G32").Value = 'Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Rana Abo-Zaid'
$ws.Range("G33").Value = 'Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Nourhan Mahmoud, Dr. Mohammad El-Tanany, Dr. Hanan Ragab, Dr. Nesma, Dr. Veronia Rafat'
$ws.Range("G34").Value = 'Dr. Amira Sobhy, Dr. Gehan Adel, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Nourhan Mahmoud, Dr. Asmaa Reda, Dr. Nahla Nagiub, Dr. Veronia Rafat, Dr. Menna tuâ€™Allah Medhat'
$ws.Range("G35").Value = 'Dr. Gehan Adel, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid'
$ws.Range("G36").Value = 'Administrator, Dr. Manar Montaser, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Shimaa Ahmad Mekki'
$ws.Range("G37").Value = 'Dr. Amira Sobhy, Dr. Gehan Adel, Dr. Manar Montaser, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid'
$ws.Range("G38").Value = 'Dr. Gehan Adel, Dr. Shimaa Ahmad Mekki, Dr. Servinaz Sayed Mohammad, Dr. Heba Mahmoud Ali, Dr. Sara Wael, Dr. Rana Abo-Zaid, Dr. Alshimaa Atef'
$ws.Range("G39").Value = 'Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Asmaa Reda'
$ws.Range("G40").Value = 'Administrator, Dr. Salma El-Gendy'
$ws.Range("G41").Value = 'Dr. Shimaa Ashraf, Dr. Safa Hany, Dr. Mariam Nour El-Din, D Wessam Atef, Dr. Omnia Mohammad'
$ws.Range("G42").Value = 'Dr. Safa Hany, Dr. Shimaa Ashraf'
$ws.Range("G43").Value = 'D Wessam Atef, Dr. Amal Awwad'
$ws.Range("G44").Value = 'Dr. Nourhan Mohammad, Dr. Amal Awwad'
$ws.Range("G45").Value = 'Dr. Marwa Mustafa, Dr. Eman M. Abo-Sakaya, Dr. Nourhan Osama, Dr. Madeha Saeed, Dr. Sarah Abdelmohsen, Dr. Yasmeena Fattoh, Dr. Basma Hamed, Dr. Esraa Mostafa, Dr. Dina Adel, Dr. Arwa Al-Sayed'
$ws.Range("G47").Value = 'D Mariam E. Mohammad, Dr. Sarah Mahdy'
$ws.Range("G50").Value = 'Dr. Nancy Abd Al-Shafy, Dr. Amr Saeed'
$ws.Range("G52").Value = 'Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Youstina Magdy, Dr. Yasmin, Dr. Monica, Dr. Maryam Ashraf, Dr. Aya Emad, Dr. Marina Atef, Dr. Salma Hassan, Dr. Remon, Dr. Neveen Nashaat'
$ws.Range("G53").Value = 'Dr. Ola Abd Al-Fattah, Dr. Youstina Magdy, Dr. Aya Emad, Dr. Abdullah El-Agrody, Dr. Marina Atef, Dr. Remon, Dr. Eman Samir Gabry'
$ws.Range("G55").Value = 'Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah, Dr. Eman Mohammad Al, Dr. Yasmin, Dr. Salma Hassan, Dr. Remon, Dr. Eman Samir Gabry, Dr. Neveen Nashaat'
$ws.Range("G56").Value = 'Dr. Aya Hanafy, Dr. Wafaa Ebida, Dr. Nardine, Dr. Abdullah El-Agrody, Dr. Salma Hassan, Dr. Remon, Dr. Eman Samir Gabry, Dr. Neveen Nashaat'
$ws.Range("G57").Value = 'Dr. Naema Gomaa, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Remon, Dr. Eman Samir Gabry, Dr. Neveen Nashaat'
